# Rebuild the LR-pair table (ECs / FAPs / sCs cross-join) with updated
# NATMI statistics, per "Natmi following Dr Hou advice".
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    ,@(2, "ECs", "ECs", 3, 1, 27.03890566666666, 81.11671699999999, 0.07096188219033728, 0.07096188219033729, 1, 0.3333333333333333, 0.1176943333333333, 0.353083, 0.05829606481791055, 0.05829606481791055, 3.182325976501222, 28.640933788511, 0.004136798483768834, 0.004136798483768835)
    ,@(3, "ECs", "FAPs", 3, 1, 27.03890566666666, 81.11671699999999, 0.07096188219033728, 0.07096188219033729, 3, 1, 1.183046666666667, 3.54914, 0.5859837360842607, 0.5859837360842608, 31.98828721926444, 287.8945849733799, 0.041582508845465, 0.04158250884546501)
    ,@(4, "ECs", "sCs", 3, 1, 27.03890566666666, 81.11671699999999, 0.07096188219033728, 0.07096188219033729, 3, 1, 0.718166, 2.154498, 0.3557201990978286, 0.3557201990978286, 19.41842272700733, 174.765804543066, 0.02524257486110344, 0.02524257486110344)
    ,@(5, "FAPs", "ECs", 3, 1, 345.566579, 1036.699737, 0.9069174311350353, 0.9069174311350354, 1, 0.3333333333333333, 0.1176943333333333, 0.353083, 0.05829606481791055, 0.05829606481791055, 40.67122813768567, 366.0410532391709, 0.05286971734994094, 0.05286971734994095)
    ,@(6, "FAPs", "FAPs", 3, 1, 345.566579, 1036.699737, 0.9069174311350353, 0.9069174311350354, 3, 1, 1.183046666666667, 3.54914, 0.5859837360842607, 0.5859837360842608, 408.8213893973533, 3679.39250457618, 0.5314388646164482, 0.5314388646164484)
    ,@(7, "FAPs", "sCs", 3, 1, 345.566579, 1036.699737, 0.9069174311350353, 0.9069174311350354, 3, 1, 0.718166, 2.154498, 0.3557201990978286, 0.3557201990978286, 248.174167774114, 2233.567509967026, 0.3226088491686461, 0.3226088491686461)
    ,@(8, "sCs", "ECs", 3, 1, 8.428738666666668, 25.286216, 0.0221206866746274, 0.02212068667462741, 1, 0.3333333333333333, 0.1176943333333333, 0.353083, 0.05829606481791055, 0.05829606481791055, 0.9920147782142223, 8.928133003928, 0.001289548984200769, 0.001289548984200769)
    ,@(9, "sCs", "FAPs", 3, 1, 8.428738666666668, 25.286216, 0.0221206866746274, 0.02212068667462741, 3, 1, 1.183046666666667, 3.54914, 0.5859837360842607, 0.5859837360842608, 9.971591183804446, 89.74432065424001, 0.01296236262234748, 0.01296236262234749)
    ,@(10, "sCs", "sCs", 3, 1, 8.428738666666668, 25.286216, 0.0221206866746274, 0.02212068667462741, 3, 1, 0.718166, 2.154498, 0.3557201990978286, 0.3557201990978286, 6.053233533285334, 54.479101799568, 0.007868775068079144, 0.007868775068079146)
)

foreach ($row in $rows) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value = $row[1]      # A: Sending cluster
    $ws.Cells.Item($r, 2).Value = "Fn1"        # B: Ligand symbol
    $ws.Cells.Item($r, 3).Value = "Tnfrsf11b"  # C: Receptor symbol
    $ws.Cells.Item($r, 4).Value = $row[2]      # D: Target cluster
    for ($i = 0; $i -lt 16; $i++) {
        $ws.Cells.Item($r, 5 + $i).Value = $row[3 + $i]
    }
}